$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new blank columns at D:E, shifting the existing D:M quarter data to F:M.
$ws.Range("D1:E1").EntireColumn.Insert()

# 2. Copy the formatting (number format / font / style) from column F (the old column D,
#    now shifted right by two) onto the two new columns D and E for the whole data block.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("F7:F102").Copy()
$ws.Range("E7:E102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# 3. Populate the two new quarter columns (D = Q4'18 12/31/2018, E = Q3'18 9/30/2018)
#    with their reported figures. $null entries are left blank (matching the rows that
#    were already blank across the rest of the table).
$data = @(
    @(7, 43465, 43373),
    @(8, 174400, 136700),
    @(9, 92000, 67000),
    @(10, 82400, 69700),
    @(11, $null, $null),
    @(12, "NA", "NA"),
    @(13, 0, 0),
    @(14, 0, 0),
    @(15, 9300, 8900),
    @(16, $null, $null),
    @(17, 121600, 93100),
    @(18, 52800, 43600),
    @(19, $null, $null),
    @(20, -3700, 2000),
    @(21, 58500, 54500),
    @(22, 15000, 17900),
    @(23, 34100, 27700),
    @(24, 5300, 7200),
    @(25, 0, 0),
    @(26, 28800, 20500),
    @(27, 24800, 17500),
    @(28, 0, 0),
    @(29, "NA", "NA"),
    @(30, 0, 0),
    @(31, 0, 0),
    @(32, 3700, -2000),
    @(33, 24800, 17500),
    @(34, 0, 0),
    @(35, 24800, 17500),
    @(38, 43465, 43373),
    @(39, $null, $null),
    @(40, $null, $null),
    @(41, 289400, 572700),
    @(42, 0, 0),
    @(43, 58300, 43500),
    @(44, 5100, 6100),
    @(45, 35700, 43100),
    @(46, 388500, 665500),
    @(47, 0, 0),
    @(48, 114400, 97200),
    @(49, 433800, 437400),
    @(50, 0, 0),
    @(51, 0, 0),
    @(52, 416700, 421300),
    @(53, 0, 0),
    @(54, 1353400, 1621400),
    @(55, $null, $null),
    @(56, $null, $null),
    @(57, 30400, 23400),
    @(58, 12000, 12000),
    @(59, 89000, 89600),
    @(60, 131400, 125000),
    @(61, 1160100, 1161700),
    @(62, 444700, 444900),
    @(63, 0, 0),
    @(64, 0, 0),
    @(65, 0, 0),
    @(66, 1728000, 1722900),
    @(67, $null, $null),
    @(68, 0, 0),
    @(69, 0, 0),
    @(70, 0, 0),
    @(71, 0, 0),
    @(72, -394400, -119000),
    @(73, 0, 0),
    @(74, 0, 0),
    @(75, 0, 0),
    @(76, -374600, -101500),
    @(77, 0, 0),
    @(80, 43465, 43373),
    @(81, 24800, 17500),
    @(82, $null, $null),
    @(83, 9400, 9000),
    @(84, 0, 0),
    @(85, 0, 0),
    @(86, 0, 0),
    @(87, 0, 0),
    @(88, 0, 0),
    @(89, 40500, 65400),
    @(90, $null, $null),
    @(91, -22300, -10500),
    @(92, 0, 0),
    @(93, 0, 0),
    @(94, -22300, -27700),
    @(95, $null, $null),
    @(96, 0, 0),
    @(97, 0, 0),
    @(98, 0, 0),
    @(99, 0, 0),
    @(100, -306200, 422900),
    @(101, -600, 200),
    @(102, -288500, 460900)
)

foreach ($item in $data) {
    $row = $item[0]
    $dVal = $item[1]
    $eVal = $item[2]
    if ($null -ne $dVal) {
        $ws.Cells.Item($row, 4).Value = $dVal
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item($row, 5).Value = $eVal
    }
}

# 4. Small restatements to the already-existing Sep-2017 (now column H) figures that came
#    along with this data refresh.
$ws.Cells.Item(20, 8).Value = 316900
$ws.Cells.Item(21, 8).Value = 367000
$ws.Cells.Item(32, 8).Value = -316900

Write-Host "Edit complete"
